# -----------------------------------------------------------------------
# Commit: "Sun, Mar 29, 2020  9:08:03 PM"
#
# What changed in the OOXML:
#   1. The table on slide 5 now points at a different <a:tableStyleId>
#      ({16BD1352-3DD4-421B-ADA6-5F43E8FA5852} -> {9E5C85D8-981E-4DCC-
#      9E06-B54D89FAC3D9}).
#   2. The presentation's theme ("Integral" / Red Violet colour scheme)
#      is replaced by the plain "Office Theme" colour scheme - i.e. the
#      author picked a different Design/colour theme for the deck.  The
#      font scheme and format scheme (fills/lines/effects) are identical
#      between the two themes, only the 12 theme colours actually change.
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-point the slide-5 table at the new table style GUID --------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{9E5C85D8-981E-4DCC-9E06-B54D89FAC3D9}")
    }
}

# --- 2. Swap the active theme's colour scheme from "Integral" / Red ---
#        Violet over to the plain "Office" palette.
function Set-ThemeRGB($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Order matches the <a:clrScheme> child order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
Set-ThemeRGB $colorScheme 1  "000000"   # dk1
Set-ThemeRGB $colorScheme 2  "FFFFFF"   # lt1
Set-ThemeRGB $colorScheme 3  "44546A"   # dk2
Set-ThemeRGB $colorScheme 4  "E7E6E6"   # lt2
Set-ThemeRGB $colorScheme 5  "5B9BD5"   # accent1
Set-ThemeRGB $colorScheme 6  "ED7D31"   # accent2
Set-ThemeRGB $colorScheme 7  "A5A5A5"   # accent3
Set-ThemeRGB $colorScheme 8  "FFC000"   # accent4
Set-ThemeRGB $colorScheme 9  "4472C4"   # accent5
Set-ThemeRGB $colorScheme 10 "70AD47"   # accent6
Set-ThemeRGB $colorScheme 11 "0563C1"   # hlink
Set-ThemeRGB $colorScheme 12 "954F72"   # folHlink
